# Add a new Job Posting row (JD_007) to the LinkedIn job postings sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 8

$ws.Cells.Item($newRow, 1).Value = "JD_007"
$ws.Cells.Item($newRow, 2).Value = "Senior Devops Developer"
$ws.Cells.Item($newRow, 3).Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Cells.Item($newRow, 4).Value = 2
$ws.Cells.Item($newRow, 5).Value = 3

# Keep the new row's height consistent with the rest of the table (avoid a
# stray custom row height from the multi-line Job_Description entry).
$ws.Rows.Item($newRow).AutoFit()
